# Mockup(Admin User).xlsx - add new params to "Table Info" sheet and
# switch the active sheet/selection to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table Info")

# The sheet grows from 6 to 9 columns (two new rows of header/data).
# Force "Text" number format first so the numeric-looking values below
# (years, ids, module numbers) are stored as shared strings (t="s"),
# exactly like the surrounding non-numeric cells, instead of being
# auto-coerced to numbers.
$ws.Range("A1:I2").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "EduYear"
$ws.Range("B1").Value = "Semester"
$ws.Range("C1").Value = "DepartmentId"
$ws.Range("D1").Value = "SpecialityId"
$ws.Range("E1").Value = "DisciplineVariantID"
$ws.Range("F1").Value = "ModuleVariantID"
$ws.Range("G1").Value = "ModuleNum"
$ws.Range("H1").Value = "NameDiscipline"
$ws.Range("I1").Value = "NameModule"

# Data row
$ws.Range("A2").Value = "2014"
$ws.Range("B2").Value = "10"
$ws.Range("C2").Value = "43"
$ws.Range("D2").Value = "6"
$ws.Range("E2").Value = "12010000137215"
$ws.Range("F2").Value = "12010000137217"
$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "Хірургія, дитяча хірургія"
$ws.Range("I2").Value = "Торакальна, серцево-судинна, ендокринна хірургія"

# Make "Table Info" the active sheet/tab with B5 selected (was "Exam
# grades" before).
$ws.Activate()
$ws.Range("B5").Select()
